$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1503.6
$ws.Range("J19").Value = 1652.2858
$ws.Range("L19").Value = 1652.2858
$ws.Range("N19").Value = -2002.2858
$ws.Range("H33").Value = 93.22222
$ws.Range("J33").Value = 122.5
$ws.Range("L33").Value = 122.5
$ws.Range("N33").Value = -580.5
$ws.Range("H64").Value = 2878.75
$ws.Range("I64").Value = 2800
$ws.Range("K64").Value = 2800
$ws.Range("M64").Value = -2552
$ws.Range("H67").Value = 2878.75
$ws.Range("I67").Value = 2800
$ws.Range("K67").Value = 2800
$ws.Range("M67").Value = -1942
$ws.Range("H125").Value = 1800
$ws.Range("J125").Value = 1700
$ws.Range("L125").Value = 15300
$ws.Range("N125").Value = -20220
$ws.Range("H132").Value = 932.19446
$ws.Range("I132").Value = 886.34375
$ws.Range("K132").Value = 2659.03125
$ws.Range("M132").Value = -129.03125
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()
$ws.Range("H135").Value = 538.8461
$ws.Range("I135").Value = 417.85715
$ws.Range("J135").Value = 680
$ws.Range("K135").Value = 3760.71435
$ws.Range("L135").Value = 6120
$ws.Range("M135").Value = -1225.71435
$ws.Range("N135").Value = -11190
$ws.Range("H138").Value = 2363.261
$ws.Range("I138").Value = 2634.3125
$ws.Range("J138").Value = 2218.7
$ws.Range("K138").Value = 7902.9375
$ws.Range("L138").Value = 6656.099999999999
$ws.Range("M138").Value = -2762.9375
$ws.Range("N138").Value = -16936.1
$ws.Range("H141").Value = 1002211.44
$ws.Range("I141").Value = 1335201.1
$ws.Range("J141").Value = 3242.5715
$ws.Range("K141").Value = 4005603.3
$ws.Range("L141").Value = 9727.7145
$ws.Range("M141").Value = -4000423.3
$ws.Range("N141").Value = -20087.7145

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3410.7083
$ws.Range("I32").Value = 2553.389
$ws.Range("K32").Value = 2553.389
$ws.Range("M32").Value = -2266.389
$ws.Range("H38").Value = 0
$ws.Range("I38").Value = 0
$ws.Range("K38").Value = 0
$ws.Range("M38").ClearContents()
$ws.Range("H45").Value = 3085.6667
$ws.Range("I45").Value = 3499.8333
$ws.Range("K45").Value = 3499.8333
$ws.Range("M45").Value = -3122.8333
$ws.Range("H74").Value = 1695.3572
$ws.Range("I74").Value = 680.1111
$ws.Range("J74").Value = 3522.8
$ws.Range("K74").Value = 680.1111
$ws.Range("L74").Value = 3522.8
$ws.Range("M74").Value = 193.8889
$ws.Range("N74").Value = -5270.8
$ws.Range("H77").Value = 1695.3572
$ws.Range("I77").Value = 680.1111
$ws.Range("J77").Value = 3522.8
$ws.Range("K77").Value = 3400.5555
$ws.Range("L77").Value = 17614
$ws.Range("M77").Value = 967.4445000000001
$ws.Range("N77").Value = -26350
$ws.Range("H104").Value = 32749.666
$ws.Range("J104").Value = 32749.666
$ws.Range("L104").Value = 32749.666
$ws.Range("N104").Value = -39737.666
$ws.Range("H132").Value = 2174.72
$ws.Range("I132").Value = 1939.6342
$ws.Range("K132").Value = 5818.902599999999
$ws.Range("M132").Value = -3288.902599999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H132").Value = 68570.71
$ws.Range("J132").Value = 68570.71
$ws.Range("L132").Value = 68570.71
$ws.Range("N132").Value = -78690.71

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1245.4634
$ws.Range("I31").Value = 773.625
$ws.Range("J31").Value = 1911.5883
$ws.Range("K31").Value = 773.625
$ws.Range("L31").Value = 1911.5883
$ws.Range("M31").Value = -478.625
$ws.Range("N31").Value = -2501.5883
$ws.Range("H34").Value = 1245.4634
$ws.Range("I34").Value = 773.625
$ws.Range("J34").Value = 1911.5883
$ws.Range("K34").Value = 773.625
$ws.Range("L34").Value = 1911.5883
$ws.Range("M34").Value = -571.625
$ws.Range("N34").Value = -2315.5883
$ws.Range("H62").Value = 1999
$ws.Range("I62").Value = 1999
$ws.Range("K62").Value = 1999
$ws.Range("M62").Value = -1375
$ws.Range("H65").Value = 1999
$ws.Range("I65").Value = 1999
$ws.Range("K65").Value = 9995
$ws.Range("M65").Value = -6875
$ws.Range("H105").Value = 1433.3334
$ws.Range("I105").Value = 1540
$ws.Range("J105").Value = 900
$ws.Range("K105").Value = 1540
$ws.Range("L105").Value = 900
$ws.Range("M105").Value = 207
$ws.Range("N105").Value = -4394
$ws.Range("H132").Value = 2254.3462
$ws.Range("I132").Value = 1535.2222
$ws.Range("K132").Value = 4605.6666
$ws.Range("M132").Value = -2075.6666
$ws.Range("H134").Value = 1732.5
$ws.Range("I134").Value = 871.2222
$ws.Range("J134").Value = 4316.3335
$ws.Range("K134").Value = 2613.6666
$ws.Range("L134").Value = 12949.0005
$ws.Range("M134").Value = -78.66660000000002
$ws.Range("N134").Value = -18019.0005

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H31").Value = 500
$ws.Range("J31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("N31").ClearContents()
$ws.Range("H92").Value = 803.8182
$ws.Range("I92").Value = 300
$ws.Range("K92").Value = 900
$ws.Range("M92").Value = 348
$ws.Range("H131").Value = 5962979.5
$ws.Range("J131").Value = 11106.5
$ws.Range("L131").Value = 33319.5
$ws.Range("N131").Value = -43399.5
$ws.Range("H140").Value = 3742.9375
$ws.Range("I140").Value = 668
$ws.Range("J140").Value = 5140.636
$ws.Range("K140").Value = 2004
$ws.Range("L140").Value = 15421.908
$ws.Range("M140").Value = 3176
$ws.Range("N140").Value = -25781.908

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 1953231.1
$ws.Range("I126").Value = 2927186.2
$ws.Range("K126").Value = 8781558.600000001
$ws.Range("M126").Value = -8779088.600000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3876.158
$ws.Range("I7").Value = 2512.7144
$ws.Range("J7").Value = 4671.5
$ws.Range("K7").Value = 2512.7144
$ws.Range("L7").Value = 4671.5
$ws.Range("M7").Value = -2400.7144
$ws.Range("N7").Value = -4895.5
$ws.Range("H40").Value = 7298.9165
$ws.Range("I40").Value = 3448.75
$ws.Range("K40").Value = 3448.75
$ws.Range("M40").Value = -3312.75
$ws.Range("H126").Value = 3876.158
$ws.Range("I126").Value = 2512.7144
$ws.Range("J126").Value = 4671.5
$ws.Range("K126").Value = 7538.1432
$ws.Range("L126").Value = 14014.5
$ws.Range("M126").Value = -5068.1432
$ws.Range("N126").Value = -18954.5
$ws.Range("H132").Value = 2279.4194
$ws.Range("I132").Value = 898.7143
$ws.Range("J132").Value = 2682.125
$ws.Range("K132").Value = 2696.1429
$ws.Range("L132").Value = 8046.375
$ws.Range("M132").Value = -166.1428999999998
$ws.Range("N132").Value = -13106.375

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("N104").ClearContents()
$ws.Range("H126").Value = 4348.68
$ws.Range("I126").Value = 3513
$ws.Range("J126").Value = 6995
$ws.Range("K126").Value = 10539
$ws.Range("L126").Value = 20985
$ws.Range("M126").Value = -8069
$ws.Range("N126").Value = -25925
$ws.Range("H132").Value = 2419.6428
$ws.Range("I132").Value = 1944.2354
$ws.Range("J132").Value = 3154.3635
$ws.Range("K132").Value = 5832.706200000001
$ws.Range("L132").Value = 9463.0905
$ws.Range("M132").Value = -3302.706200000001
$ws.Range("N132").Value = -14523.0905
$ws.Range("H136").Value = 25256946
$ws.Range("I136").Value = 55560480
$ws.Range("J136").Value = 4000
$ws.Range("K136").Value = 166681440
$ws.Range("L136").Value = 12000
$ws.Range("M136").Value = -166678890
$ws.Range("N136").Value = -17100
